$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.252.68'
$ws.Range("E2").Value = '  +1.60%  '
$ws.Range("D3").Value = '2.655.94'
$ws.Range("E3").Value = '  +1.95%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.86'
$ws.Range("E5").Value = '  +2.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.00'
$ws.Range("E6").Value = '  +1.50%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.572'
$ws.Range("E8").Value = '  +1.47%  '
$ws.Range("D9").Value = '2.666.43'
$ws.Range("E9").Value = '  +2.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.34'
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("E11").Value = '  +1.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.338'
$ws.Range("E12").Value = '  +1.37%  '
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("D14").Value = '3.118.64'
$ws.Range("E14").Value = '  +2.13%  '
$ws.Range("D15").Value = '59.269.09'
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.01'
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("E17").Value = '  +1.32%  '
$ws.Range("D18").Value = '2.660.37'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '353.16'
$ws.Range("E19").Value = '  +4.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.51'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.43'
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("E22").Value = '  +3.20%  '
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.98'
$ws.Range("E24").Value = '  +2.97%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.165'
$ws.Range("E26").Value = '  +4.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.14'
$ws.Range("E29").Value = '  +1.90%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.59'
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.99'
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.25'
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.955'
$ws.Range("E35").Value = '  -7.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.08'
$ws.Range("E36").Value = '  +3.63%  '
$ws.Range("E37").Value = '  +5.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.865'
$ws.Range("E38").Value = '  +1.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.68'
$ws.Range("E39").Value = '  +1.84%  '
$ws.Range("E40").Value = '  +3.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.68'
$ws.Range("E41").Value = '  +1.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '281.20'
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("E43").Value = '  +0.73%  '
$ws.Range("E44").Value = '  +0.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.86'
$ws.Range("E45").Value = '  +3.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.606'
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").Value = '2.123.24'
$ws.Range("E47").Value = '  +8.94%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.85'
$ws.Range("E48").Value = '  +3.30%  '
$ws.Range("B49").Value = 'Hedera'
$ws.Range("C49").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0530'
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("E50").Value = '  +1.55%  '
